$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (new crime data collected: volume bump + report week roll) ---
$ws.Range("A8").Value = "Volume 32   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  3/31/2025  Through  4/6/2025"

# --- Weekly crime-complaint table (rows 14-28, 31) ---
# Row 14
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G14").Value = 1
$ws.Range("G14").NumberFormat = '#,##0'
$ws.Range("H14").Value = -100
$ws.Range("H14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 50

# Row 15
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("N15").Value = -70

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = -18.918918918918
$ws.Range("M16").Value = -57.746478873239
$ws.Range("N16").Value = -88.235294117647

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 61
$ws.Range("J17").Value = 72
$ws.Range("K17").Value = -15.277777777777
$ws.Range("L17").Value = -1.612903225806
$ws.Range("M17").Value = 190.47619047619
$ws.Range("N17").Value = 1.666666666666

# Row 18
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 28.571428571428
$ws.Range("F18").Value = 29
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 81.25
$ws.Range("I18").Value = 72
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = 16.129032258064
$ws.Range("L18").Value = 2.857142857142
$ws.Range("M18").Value = 4.347826086956
$ws.Range("N18").Value = -82.178217821782

# Row 19
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -6.666666666666
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -8.333333333333
$ws.Range("I19").Value = 118
$ws.Range("J19").Value = 162
$ws.Range("K19").Value = -27.160493827160
$ws.Range("L19").Value = -28.915662650602
$ws.Range("M19").Value = -4.838709677419
$ws.Range("N19").Value = -14.492753623188

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 150
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = 8.333333333333
$ws.Range("I20").Value = 89
$ws.Range("J20").Value = 75
$ws.Range("K20").Value = 18.666666666666
$ws.Range("L20").Value = 30.882352941176
$ws.Range("M20").Value = 78
$ws.Range("N20").Value = -93.019607843137

# Row 21
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = 10.714285714285
$ws.Range("F21").Value = 125
$ws.Range("G21").Value = 116
$ws.Range("H21").Value = 7.758620689655
$ws.Range("I21").Value = 376
$ws.Range("J21").Value = 410
$ws.Range("K21").Value = -8.292682926829
$ws.Range("L21").Value = -8.292682926829
$ws.Range("M21").Value = 11.242603550295
$ws.Range("N21").Value = -82.462686567164

# Row 22
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("L14").Copy($ws.Range("E22"))
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -100
$ws.Range("L22").Value = -54.545454545454

# Row 23
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("L14").Copy($ws.Range("E23"))
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = 150

# Row 24
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 4.761904761904
$ws.Range("F24").Value = 80
$ws.Range("G24").Value = 97
$ws.Range("H24").Value = -17.525773195876
$ws.Range("I24").Value = 262
$ws.Range("J24").Value = 315
$ws.Range("K24").Value = -16.825396825396
$ws.Range("L24").Value = -33.163265306122
$ws.Range("M24").Value = 34.358974358974

# Row 25
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 30
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = -45.454545454545
$ws.Range("I25").Value = 95
$ws.Range("J25").Value = 150
$ws.Range("K25").Value = -36.666666666666
$ws.Range("L25").Value = -26.923076923076

# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 31.034482758620
$ws.Range("I26").Value = 149
$ws.Range("J26").Value = 110
$ws.Range("K26").Value = 35.454545454545
$ws.Range("L26").Value = 15.503875968992
$ws.Range("M26").Value = 46.078431372549

# Row 27
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("L27").Value = -60

# Row 28
$ws.Range("C28").Value = 3
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 18
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = 63.636363636363
$ws.Range("L28").Value = 63.636363636363

# Row 31
$ws.Range("L31").Value = -50

